$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-05-12 Sunday" "2024-05-13 Monday"

Replace-Text "576×5=" "135×2="
Replace-Text "443×3=" "579×2="
Replace-Text "870×6=" "902×3="
Replace-Text "492×3=" "726×2="
Replace-Text "909×2=" "455×7="
Replace-Text "680×4=" "419×9="
Replace-Text "149×9=" "193×2="
Replace-Text "181×8=" "132×2="
Replace-Text "628×7=" "232×6="
Replace-Text "228×4=" "717×5="
Replace-Text "691×6=" "826×2="
Replace-Text "517×2=" "682×2="
Replace-Text "375×2=" "185×8="
Replace-Text "491×6=" "623×4="
Replace-Text "644×8=" "976×4="
Replace-Text "237×4=" "714×6="
Replace-Text "895×2=" "486×6="
Replace-Text "291×9=" "446×3="
Replace-Text "918×2=" "838×5="
Replace-Text "448×5=" "668×9="
Replace-Text "684×6=" "958×7="
Replace-Text "357×3=" "791×5="
Replace-Text "934×9=" "898×6="
Replace-Text "224×5=" "830×7="
Replace-Text "473×3=" "540×7="
